$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - 09/03/2018
$ws.Range("A13").Value = 43168
$ws.Range("B13").Value = "Discussion de sur l'implémentation du undo-redo + établissement d'un diagramme de classe du model"
$ws.Range("C13").Value = 4
$ws.Rows.Item(13).RowHeight = 30

# Row 14 - 12/03/2018
$ws.Range("A14").Value = 43171
$ws.Range("B14").Value = "Réunion de mise en commun du travail + discussion du undo-redo"
$ws.Range("C14").Value = 1.5
$ws.Rows.Item(14).RowHeight = 30

# Row 15 - 16/03/2018
$ws.Range("A15").Value = 43175
$ws.Range("B15").Value = "Recherche sur l'implementation de la sauvegarde"
$ws.Range("C15").Value = 2

# Row 16 - 18/03/2018
$ws.Range("A16").Value = 43177
$ws.Range("B16").Value = "test d'implémentation du système de sauvegrade par sérialisation"
$ws.Range("C16").Value = 5
$ws.Rows.Item(16).RowHeight = 30

# Update the active selection to match the saved view state
[void]$ws.Range("B15").Select()
